# Daily update: append the next day's row of win counts to the bottom
# of the tracking sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 50
$newRow = $lastRow + 1

# The last populated row in the "Day" column carries a special
# "last row" date style. Move that style down to the newly appended
# row, and restore the previous last row to the normal per-row date
# style used throughout the rest of the column.
$ws.Cells.Item($newRow, 1).NumberFormat = $ws.Cells.Item($lastRow, 1).NumberFormat
$ws.Cells.Item($lastRow, 1).NumberFormat = $ws.Cells.Item($lastRow - 1, 1).NumberFormat

# Append today's data.
$ws.Cells.Item($newRow, 1).Value = 45791
$ws.Cells.Item($newRow, 2).Value = 212
$ws.Cells.Item($newRow, 3).Value = 217
$ws.Cells.Item($newRow, 4).Value = 219
